$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2075.5715
$ws.Range("I32").Value = 994.5
$ws.Range("K32").Value = 994.5
$ws.Range("M32").Value = -668.5
$ws.Range("H43").Value = 1464.8182
$ws.Range("I43").Value = 1500
$ws.Range("J43").Value = 1461.3
$ws.Range("K43").Value = 1500
$ws.Range("L43").Value = 1461.3
$ws.Range("M43").Value = -1431
$ws.Range("N43").Value = -1599.3
$ws.Range("H55").Value = 184.4
$ws.Range("I55").Value = 200
$ws.Range("K55").Value = 200
$ws.Range("M55").Value = 14
$ws.Range("H69").Value = 3777.4
$ws.Range("I69").Value = 2962.3333
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 8886.999899999999
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -8012.999899999999
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 3777.4
$ws.Range("I72").Value = 2962.3333
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 26660.9997
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -22292.9997
$ws.Range("N72").Value = -53736
$ws.Range("H98").Value = 2575.7058
$ws.Range("J98").Value = 2131.3333
$ws.Range("L98").Value = 2131.3333
$ws.Range("N98").Value = -5127.3333
$ws.Range("H112").Value = 1460.2703
$ws.Range("J112").Value = 1473.0555
$ws.Range("L112").Value = 4419.166499999999
$ws.Range("N112").Value = -6635.166499999999
$ws.Range("H122").Value = 2575.7058
$ws.Range("J122").Value = 2131.3333
$ws.Range("L122").Value = 6393.999899999999
$ws.Range("N122").Value = -11293.9999
$ws.Range("H137").Value = 1814.8334
$ws.Range("I137").Value = 997.25
$ws.Range("K137").Value = 2991.75
$ws.Range("M137").Value = -441.75
$ws.Range("H138").Value = 2782.4465
$ws.Range("J138").Value = 2248.639
$ws.Range("L138").Value = 6745.917
$ws.Range("N138").Value = -17025.917
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4093.9033
$ws.Range("I32").Value = 3252.158
$ws.Range("J32").Value = 13689.8
$ws.Range("K32").Value = 3252.158
$ws.Range("L32").Value = 13689.8
$ws.Range("M32").Value = -2965.158
$ws.Range("N32").Value = -14263.8
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = ""
$ws.Range("H61").Value = 7595.353
$ws.Range("I61").Value = 8055.5386
$ws.Range("J61").Value = 6099.75
$ws.Range("K61").Value = 8055.5386
$ws.Range("L61").Value = 6099.75
$ws.Range("M61").Value = -7843.5386
$ws.Range("N61").Value = -6523.75
$ws.Range("H102").Value = 1344.3077
$ws.Range("I102").Value = 1157.4
$ws.Range("K102").Value = 1157.4
$ws.Range("M102").Value = 464.5999999999999
$ws.Range("H123").Value = 68000
$ws.Range("J123").Value = 68000
$ws.Range("L123").Value = 68000
$ws.Range("N123").Value = -77800
$ws.Range("H136").Value = 7595.353
$ws.Range("I136").Value = 8055.5386
$ws.Range("J136").Value = 6099.75
$ws.Range("K136").Value = 24166.6158
$ws.Range("L136").Value = 18299.25
$ws.Range("M136").Value = -21616.6158
$ws.Range("N136").Value = -23399.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620
$ws.Range("H94").Value = 1484.5
$ws.Range("I94").Value = 477.25
$ws.Range("J94").Value = 3499
$ws.Range("K94").Value = 477.25
$ws.Range("L94").Value = 3499
$ws.Range("M94").Value = -26.25
$ws.Range("N94").Value = -4401
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2680.75
$ws.Range("J31").Value = 2922.8
$ws.Range("L31").Value = 2922.8
$ws.Range("N31").Value = -3512.8
$ws.Range("H34").Value = 2680.75
$ws.Range("J34").Value = 2922.8
$ws.Range("L34").Value = 2922.8
$ws.Range("N34").Value = -3326.8
$ws.Range("H58").Value = 2418151.8
$ws.Range("I58").Value = 3954776.8
$ws.Range("K58").Value = 3954776.8
$ws.Range("M58").Value = -3954573.8
$ws.Range("H86").Value = 3210.25
$ws.Range("I86").Value = 3229.3635
$ws.Range("K86").Value = 3229.3635
$ws.Range("M86").Value = -2106.3635
$ws.Range("H89").Value = 3210.25
$ws.Range("I89").Value = 3229.3635
$ws.Range("K89").Value = 16146.8175
$ws.Range("M89").Value = -10530.8175
$ws.Range("H105").Value = 847.6923
$ws.Range("I105").Value = 827.63635
$ws.Range("K105").Value = 827.63635
$ws.Range("M105").Value = 919.36365
$ws.Range("H136").Value = 2418151.8
$ws.Range("I136").Value = 3954776.8
$ws.Range("K136").Value = 11864330.4
$ws.Range("M136").Value = -11861780.4
$ws.Range("H141").Value = 62249.25
$ws.Range("J141").Value = 58999
$ws.Range("L141").Value = 58999
$ws.Range("N141").Value = -69359
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 109.75
$ws.Range("I40").Value = 96.333336
$ws.Range("J40").Value = 150
$ws.Range("K40").Value = 385.333344
$ws.Range("L40").Value = 600
$ws.Range("M40").Value = -316.333344
$ws.Range("N40").Value = -738
$ws.Range("H81").Value = 1512.1428
$ws.Range("J81").Value = 2301.6667
$ws.Range("L81").Value = 6905.000100000001
$ws.Range("N81").Value = -9151.000100000001
$ws.Range("H84").Value = 1512.1428
$ws.Range("J84").Value = 2301.6667
$ws.Range("L84").Value = 20715.0003
$ws.Range("N84").Value = -31947.0003
$ws.Range("H107").Value = 486
$ws.Range("I107").Value = 416.66666
$ws.Range("J107").Value = 504.9091
$ws.Range("K107").Value = 1249.99998
$ws.Range("L107").Value = 1514.7273
$ws.Range("M107").Value = 670.0000199999999
$ws.Range("N107").Value = -5354.7273
$ws.Range("H110").Value = 1999.6666
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H131").Value = 15983.805
$ws.Range("J131").Value = 17851.854
$ws.Range("L131").Value = 53555.562
$ws.Range("N131").Value = -63635.562
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3491900
$ws.Range("I7").Value = 5375000
$ws.Range("J7").Value = 667249.75
$ws.Range("K7").Value = 5375000
$ws.Range("L7").Value = 667249.75
$ws.Range("M7").Value = -5374888
$ws.Range("N7").Value = -667473.75
$ws.Range("H8").Value = 3491900
$ws.Range("I8").Value = 5375000
$ws.Range("J8").Value = 667249.75
$ws.Range("K8").Value = 5375000
$ws.Range("L8").Value = 667249.75
$ws.Range("M8").Value = -5374861
$ws.Range("N8").Value = -667527.75
$ws.Range("H11").Value = 6046900.5
$ws.Range("I11").Value = 6784179
$ws.Range("J11").Value = 4045716.5
$ws.Range("K11").Value = 6784179
$ws.Range("L11").Value = 4045716.5
$ws.Range("M11").Value = -6784040
$ws.Range("N11").Value = -4045994.5
$ws.Range("H12").Value = 5427647.5
$ws.Range("I12").Value = 6500000
$ws.Range("J12").Value = 2854002.2
$ws.Range("K12").Value = 6500000
$ws.Range("L12").Value = 2854002.2
$ws.Range("M12").Value = -6499860
$ws.Range("N12").Value = -2854282.2
$ws.Range("H80").Value = 3156.6667
$ws.Range("J80").Value = 2490
$ws.Range("L80").Value = 2490
$ws.Range("N80").Value = -4486
$ws.Range("H83").Value = 3156.6667
$ws.Range("J83").Value = 2490
$ws.Range("L83").Value = 12450
$ws.Range("N83").Value = -22434
$ws.Range("H136").Value = 13105.714
$ws.Range("J136").Value = 13105.714
$ws.Range("L136").Value = 39317.142
$ws.Range("N136").Value = -44417.142
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1673.3684
$ws.Range("I82").Value = 1522.125
$ws.Range("J82").Value = 1783.3636
$ws.Range("K82").Value = 1522.125
$ws.Range("L82").Value = 1783.3636
$ws.Range("M82").Value = -1161.125
$ws.Range("N82").Value = -2505.3636
$ws.Range("H85").Value = 1673.3684
$ws.Range("I85").Value = 1522.125
$ws.Range("J85").Value = 1783.3636
$ws.Range("K85").Value = 1522.125
$ws.Range("L85").Value = 1783.3636
$ws.Range("M85").Value = -274.125
$ws.Range("N85").Value = -4279.3636
$ws.Range("H136").Value = 2420.8
$ws.Range("I136").Value = 2050.75
$ws.Range("K136").Value = 6152.25
$ws.Range("M136").Value = -3602.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = ""
$ws.Range("H12").Value = 65005.25
$ws.Range("J12").Value = 76673.664
$ws.Range("L12").Value = 76673.664
$ws.Range("N12").Value = -76957.664
$ws.Range("H14").Value = 4998.5
$ws.Range("J14").Value = 4998.5
$ws.Range("L14").Value = 4998.5
$ws.Range("N14").Value = -5334.5
$ws.Range("H81").Value = 942.8889
$ws.Range("I81").Value = 967
$ws.Range("K81").Value = 1934
$ws.Range("M81").Value = -873
$ws.Range("H84").Value = 942.8889
$ws.Range("I84").Value = 967
$ws.Range("K84").Value = 9670
$ws.Range("M84").Value = -4366
$ws.Range("H113").Value = 535.2105
$ws.Range("I113").Value = 316.6
$ws.Range("K113").Value = 949.8000000000001
$ws.Range("M113").Value = 1220.2
$ws.Range("H132").Value = 2580.7778
$ws.Range("I132").Value = 1657.6923
$ws.Range("K132").Value = 4973.0769
$ws.Range("M132").Value = -2443.0769
$ws.Range("H136").Value = 21368960
$ws.Range("I136").Value = 32680752
$ws.Range("J136").Value = 2239.111
$ws.Range("K136").Value = 98042256
$ws.Range("L136").Value = 6717.333
$ws.Range("M136").Value = -98039706
$ws.Range("N136").Value = -11817.333
